$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.994.32"
$ws.Range("E2").Value = "  -3.82%  "
$ws.Range("D3").Value = "1.867.71"
$ws.Range("E3").Value = "  -3.02%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'318.64"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4328"
$ws.Range("E7").Value = "  -6.17%  "
$ws.Range("D8").Value = "'0.3713"
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("D9").Value = "'0.07403"
$ws.Range("E9").Value = "  -4.70%  "
$ws.Range("D10").Value = "'0.9293"
$ws.Range("E10").Value = "  -5.18%  "
$ws.Range("E11").Value = "  -6.56%  "
$ws.Range("D12").Value = "1.852.63"
$ws.Range("E12").Value = "  -5.21%  "
$ws.Range("D13").Value = "'6.722"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").Value = "'5.420"
$ws.Range("E14").Value = "  -5.01%  "
$ws.Range("D15").Value = "'0.06881"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'80.47"
$ws.Range("E17").Value = "  -4.71%  "
$ws.Range("D18").Value = "'0.000009006"
$ws.Range("E18").Value = "  -5.67%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'15.73"
$ws.Range("E20").Value = "  -6.18%  "
$ws.Range("D21").Value = "27.991.94"
$ws.Range("E21").Value = "  -3.82%  "
$ws.Range("D22").Value = "'5.112"
$ws.Range("E22").Value = "  -4.51%  "
$ws.Range("D23").Value = "'10.99"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "2.199.64"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").Value = "'2.049"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").Value = "'154.03"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").Value = "'5.493"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "'112.78"
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("D30").Value = "'1.687"
$ws.Range("E30").Value = "  -8.55%  "
$ws.Range("D31").Value = "'0.08972"
$ws.Range("E31").Value = "  -4.07%  "
$ws.Range("D32").Value = "'0.8047"
$ws.Range("E32").Value = "  -6.47%  "
$ws.Range("D33").Value = "'4.766"
$ws.Range("E33").Value = "  -7.00%  "
$ws.Range("D34").Value = "'1.172"
$ws.Range("E34").Value = "  -5.97%  "
$ws.Range("D35").Value = "'2.952"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").Value = "'1.006"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "'1.119"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("D39").Value = "'0.01971"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("D40").Value = "'2.994"
$ws.Range("E40").Value = "  -6.92%  "
$ws.Range("D41").Value = "'0.5218"
$ws.Range("E41").Value = "  -5.60%  "
$ws.Range("D42").Value = "'6.986"
$ws.Range("E42").Value = "  -7.12%  "
$ws.Range("D43").Value = "'0.1680"
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("D44").Value = "'8.719"
$ws.Range("E44").Value = "  -6.86%  "
$ws.Range("D45").Value = "'0.06707"
$ws.Range("E45").Value = "  -3.22%  "
$ws.Range("D46").Value = "'0.4865"
$ws.Range("E46").Value = "  -6.66%  "
$ws.Range("D47").Value = "'10.48"
$ws.Range("E47").Value = "  -7.50%  "
$ws.Range("D48").Value = "'106.67"
$ws.Range("E48").Value = "  -3.45%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "'1.666"
$ws.Range("E50").Value = "  -6.04%  "
$ws.Range("D51").Value = "'1.866"
$ws.Range("E51").Value = "  -15.06%  "
